# Apply minor date corrections to the TrialData sheet and update the
# active cell selection, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrialData")

# Date value corrections (stored as raw Excel serial date numbers so the
# existing date number-format style on these cells is preserved).
$ws.Range("Q2").Value = 42108
$ws.Range("O3").Value = 42470
$ws.Range("Q4").Value = 42840
$ws.Range("O5").Value = 43177
$ws.Range("Q5").Value = 43193
$ws.Range("O6").Value = 43557
$ws.Range("Q6").Value = 43572
$ws.Range("Q7").Value = 43935
$ws.Range("O8").Value = 44301

# Update the selected/active cell shown when the workbook is reopened.
$ws.Range("Q9").Select()
